$wb = $excel.ActiveWorkbook

$wsSurvey = $wb.Worksheets.Item("survey")
$wsSettings = $wb.Worksheets.Item("settings")

# Update translation-key labels (order matters for shared-string table layout)
$wsSettings.Range("C1").Value = "display.title.text"
$wsSurvey.Range("C1").Value = "display.prompt.text"
$wsSurvey.Range("D1").Value = "display.hint.text"

# Row 2 on the survey sheet grows to fit the longer wrapped label text
$wsSurvey.Rows.Item(2).RowHeight = 54.35

# Update sheet selections; select settings first, survey last so survey ends
# up as the active/selected tab (mirrors the workbook-level activeTab change)
$null = $wsSettings.Range("C2").Select()
$null = $wsSurvey.Range("D2").Select()
